# Add a new "ds" / "beer" / "mouth Fresh" set of item+amount pairs to row 2
# of the active sheet ("08-07-2023"), mirroring the existing layout
# (item text columns alternating with currency amount columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2: amount for existing item in C2, formatted like the other currency cells (style index 1)
$ws.Range("D2").Value = 18800
$ws.Range("D2").NumberFormat = $ws.Range("B2").NumberFormat

# E2 / F2: "ds" item + amount, new currency-like number format (style index 2 / numFmtId 164)
$ws.Range("E2").Value = "ds"
$ws.Range("F2").Value = 11200
$ws.Range("F2").NumberFormat = "_ [$₹-4009]\ * #,##0.00_ ;_ [$₹-4009]\ * \-#,##0.00_ ;_ [$₹-4009]\ * ""-""??_ ;_ @_ "

# G2 / H2: "beer" item + amount
$ws.Range("G2").Value = "beer"
$ws.Range("H2").Value = 100
$ws.Range("H2").NumberFormat = "_ [$₹-4009]\ * #,##0.00_ ;_ [$₹-4009]\ * \-#,##0.00_ ;_ [$₹-4009]\ * ""-""??_ ;_ @_ "

# I2: "mouth Fresh" item
$ws.Range("I2").Value = "mouth Fresh"

# New amount columns get the same display width as column B (the existing amount column)
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(2).ColumnWidth
$ws.Columns.Item(6).ColumnWidth = $ws.Columns.Item(2).ColumnWidth

$ws.Range("I3").Select()
